# Fix Training Data Issue (#48)
#
# The "Date" column (BF) stored the literal text "6-5-2013-14" for every
# team row on this sheet. The NBA stats site showed game dates one day
# off from how they were actually recorded, so the training data needs
# to be corrected to the real ISO date "2014-06-05" for every row
# (BF2:BF31).
#
# Setting a date-looking string straight into .Value would make Excel
# auto-convert it into a real date serial number, which is not what we
# want here - the column must keep holding literal text. So each cell is
# switched to a text number format before the value is written, and the
# format is cleared again afterwards so the cells end up with no extra
# styling applied (matching their original unstyled state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateColumn = 58          # column BF
$firstRow = 2
$lastRow = 31
$correctDate = "2014-06-05"

$dateRange = $ws.Range($ws.Cells.Item($firstRow, $dateColumn), $ws.Cells.Item($lastRow, $dateColumn))
$dateRange.NumberFormat = "@"

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, $dateColumn).Value = $correctDate
}

$dateRange.ClearFormats()
